$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B9: convert the "Cedula" value from an inline/text string to a real numeric value
$ws.Range("B9").Value = 1000127336

# New row 10: a fresh negotiation-log entry
$ws.Range("A10").Value = "2025-10-16 12:49:16"

# B10 stays textual (like the header "Cedula" column for every other row) even
# though it looks numeric, so force text formatting before assigning it, then
# drop the number-format override so the cell keeps the sheet's default style.
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "1000274330"
$ws.Range("B10").ClearFormats()

$ws.Range("C10").Value = "Elian"
$ws.Range("D10").Value = "TARJETA DE CRÉDITO"
$ws.Range("E10").Value = "****0786"
$ws.Range("F10").Value = "PRORROGA SIN PAGO"
$ws.Range("G10").Value = "24 cuotas"
$ws.Range("H10").Value = "35.197.92.111"
$ws.Range("I10").Value = "The Dalles"
$ws.Range("J10").Value = "Oregon"
$ws.Range("K10").Value = "United States"
$ws.Range("L10").Value = "2025-10-16 12:49:16"
$ws.Range("M10").Value = "*****0786"
$ws.Range("N10").Value = "35.197.92.111"
# O10 (MensajeUsuario) and P10 (RespuestaIA) are left blank for this row,
# same as every other column in the sheet that has no content.
